$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$shp = $nm.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
Write-Output "before: [$($tr.Text)]"
$tr.Text = "4/10/2025"
Write-Output "after: [$($tr.Text)]"
